$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-03-09 Sunday" "2025-03-10 Monday"

Replace-Text "705÷4=176, 1" "312÷2=156, 0"
Replace-Text "372÷9=41, 3" "388÷2=194, 0"
Replace-Text "796÷7=113, 5" "444÷5=88, 4"
Replace-Text "289÷3=96, 1" "819÷6=136, 3"
Replace-Text "346÷9=38, 4" "837÷3=279, 0"

Replace-Text "297÷8=37, 1" "778÷8=97, 2"
Replace-Text "304÷2=152, 0" "264÷5=52, 4"
Replace-Text "607÷9=67, 4" "328÷8=41, 0"
Replace-Text "734÷7=104, 6" "741÷9=82, 3"
Replace-Text "225÷7=32, 1" "638÷8=79, 6"

Replace-Text "345÷4=86, 1" "588÷6=98, 0"
Replace-Text "921÷4=230, 1" "257÷9=28, 5"
Replace-Text "131÷9=14, 5" "261÷2=130, 1"
Replace-Text "269÷2=134, 1" "309÷4=77, 1"
Replace-Text "486÷7=69, 3" "324÷5=64, 4"

Replace-Text "969÷9=107, 6" "120÷9=13, 3"
Replace-Text "566÷4=141, 2" "344÷9=38, 2"
Replace-Text "513÷3=171, 0" "194÷2=97, 0"
Replace-Text "819÷9=91, 0" "254÷4=63, 2"
Replace-Text "565÷7=80, 5" "575÷4=143, 3"

Replace-Text "835÷8=104, 3" "132÷3=44, 0"
Replace-Text "440÷2=220, 0" "606÷7=86, 4"
Replace-Text "329÷4=82, 1" "619÷7=88, 3"
Replace-Text "914÷4=228, 2" "780÷7=111, 3"
Replace-Text "606÷6=101, 0" "190÷5=38, 0"
